$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 6 (Бусин Дмитрий): fill G6/H6, add J6, extend sum formula to include H
# ---------------------------------------------------------------------------
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 5
$ws.Range("J6").Value = 5
$ws.Range("K6").Formula = "=SUM(C6:H6)"

# ---------------------------------------------------------------------------
# Row 10 (Иванов Вячеслав): C10:G10 raised from 2 to 5 and restyled from the
# "ДЗ not done" (green) look to the plain thick-border look used by G10/H10
# ---------------------------------------------------------------------------
$ws.Range("H10").Copy()
$ws.Range("C10:G10").PasteSpecial(-4122)
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 5

# ---------------------------------------------------------------------------
# Row 15 (Катахова Марина): fill G15/H15, add J15, extend sum formula, and
# give the row the thicker "bottom of block" look (taller + thick top/bottom)
# ---------------------------------------------------------------------------
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 5
$ws.Range("J15").Value = 5
$ws.Range("K15").Formula = "=SUM(C15:H15)"
$ws.Rows(15).RowHeight = 14

# ---------------------------------------------------------------------------
# Row 16 (Кобзев Богдан): fill H16, add I16/J16 (matching the style already
# used on G16/H16), give the row the same thicker look as row 15
# ---------------------------------------------------------------------------
$ws.Range("H16").Value = 5
$ws.Range("G16").Copy()
$ws.Range("I16:J16").PasteSpecial(-4122)
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 5
$ws.Rows(16).RowHeight = 14

# ---------------------------------------------------------------------------
# Row 17 (Корюгин Андрей): only the thicker row look changes, no cell edits
# ---------------------------------------------------------------------------
$ws.Rows(17).RowHeight = 14

# ---------------------------------------------------------------------------
# Row 20 (Одинцев Константин): F20 raised from 2 to 5, restyled to match the
# plain thick-border look (like C20:E20)
# ---------------------------------------------------------------------------
$ws.Range("E20").Copy()
$ws.Range("F20").PasteSpecial(-4122)
$ws.Range("F20").Value = 5

# ---------------------------------------------------------------------------
# Restore the selection to J6, as recorded in the saved sheet view
# ---------------------------------------------------------------------------
$ws.Range("J6").Select()
